$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: bay "2-1" belongs to hangar 2, not hangar 1 as previously listed.
# Add the corrected row plus the other rows that were captured while
# reproducing/validating the issue.

# Row 3: Bay 2-1, Hangar 2, SN123, Test Customer, Rank 1
$ws.Range("A3").Value = "2-1"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "SN123"
$ws.Range("D3").Value = "Test Customer"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "https://example.com/dashboard|https://example.com/status"

# Row 4: Bay 1-1, Hangar 2, 56887, Carbo Inc., Rank 3
$ws.Range("A4").Value = "1-1"
$ws.Range("B4").Value = 2
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "56887"
$ws.Range("D4").Value = "Carbo Inc."
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = "https://powerbi.bellflight.com/reports/powerbi/Piney%20Flats/Aircraft%20Services/Part%20Visibility%20Report?rs:embed=true|https://powerbi.bellflight.com/reports/powerbi/Piney%20Flats/Aircraft%20Services/Part%20Visibility%20Report?rs:embed=true"

# Row 5: Bay 1-1, Hangar 1, 56887, hass bombn, Rank 3
$ws.Range("A5").Value = "1-1"
$ws.Range("B5").Value = 1
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "56887"
$ws.Range("D5").Value = "hass bombn"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = "https://powerbi.bellflight.com/reports/powerbi/Piney%20Flats/Aircraft%20Services/Part%20Visibility%20Report?rs:embed=true"
